$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 - new amortization case
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 3956
$ws.Range("C27").Value = 1665214
$ws.Range("D27").Value = "GALLARDO"
$ws.Range("E27").Value = "VACAFLOR"
$ws.Range("F27").Value = "CRISOSTOMO"
$ws.Range("J27").Value = "VEJEZ"
$ws.Range("K27").Value = 7431.04
$ws.Range("L27").Value = 4696.5600000000004
$ws.Range("M27").Value = 2734.48

# Row 28 - new amortization case
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = 3955
$ws.Range("C28").Value = 2323723
$ws.Range("J28").Value = "VEJEZ"
$ws.Range("K28").Value = 7503.18
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 7503.18

# Apply autofilter on F1:F26
$ws.Range("F1:F26").AutoFilter()

# Update selection
$ws.Range("O20").Select()
